# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets
# to reflect newly scraped totals (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# Row => new value, for the "展览" sheet (sheet1)
$exhibitUpdates = @{
    5  = 334
    8  = 1100
    9  = 307
    14 = 548
    15 = 8628
    17 = 2487
    18 = 257
    22 = 589
    24 = 1154
    25 = 996
    26 = 2025
    27 = 2066
    29 = 1760
    34 = 48
    40 = 246
    41 = 426
    42 = 683
}

foreach ($row in $exhibitUpdates.Keys) {
    $wsExhibit.Range("F$row").Value = $exhibitUpdates[$row]
}

# Row => new value, for the "全部类型" sheet (sheet4)
$allUpdates = @{
    5  = 334
    9  = 1100
    10 = 307
    14 = 548
    15 = 8628
    17 = 2487
    19 = 257
    23 = 589
    25 = 1154
    26 = 996
    27 = 2025
    28 = 2066
    29 = 1760
    34 = 48
    40 = 246
    41 = 426
    46 = 683
}

foreach ($row in $allUpdates.Keys) {
    $wsAll.Range("F$row").Value = $allUpdates[$row]
}
